$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python")

# Row 37 already existed in the sheet; this upload only added a D37 value
$ws.Cells.Item(37,4).Value = "completed"

# Row 38
$ws.Cells.Item(38,1).Value = 45695
$ws.Cells.Item(38,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(38,2).Value = "sql practice-rename colum,colum delete,values delete"
$ws.Cells.Item(38,4).Value = "completed"

# Row 39
$ws.Cells.Item(39,1).Value = 45696
$ws.Cells.Item(39,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(39,2).Value = "Saturday"

# Row 40
$ws.Cells.Item(40,1).Value = 45697
$ws.Cells.Item(40,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(40,2).Value = "Sunday"

# Row 41
$ws.Cells.Item(41,1).Value = 45698
$ws.Cells.Item(41,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(41,2).Value = "leave"

# Row 42
$ws.Cells.Item(42,1).Value = 45699
$ws.Cells.Item(42,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(42,2).Value = "leave"

# Row 43
$ws.Cells.Item(43,1).Value = 45700
$ws.Cells.Item(43,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(43,2).Value = "Assignment 11(append)"
$ws.Cells.Item(43,4).Value = "completed"

# Row 44
$ws.Cells.Item(44,1).Value = 45701
$ws.Cells.Item(44,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(44,2).Value = "Assignment 11(append)"
$ws.Cells.Item(44,4).Value = "completed"

# Row 45
$ws.Cells.Item(45,1).Value = 45702
$ws.Cells.Item(45,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(45,2).Value = "Assignment11(append),SQL practice"
$ws.Cells.Item(45,4).Value = "completed"

# Row 46
$ws.Cells.Item(46,1).Value = 45703
$ws.Cells.Item(46,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(46,2).Value = "Saturday"

# Row 47
$ws.Cells.Item(47,1).Value = 45704
$ws.Cells.Item(47,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(47,2).Value = "Sunday"

# Row 48
$ws.Cells.Item(48,1).Value = 45705
$ws.Cells.Item(48,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(48,2).Value = "Assignment12(merges)"
$ws.Cells.Item(48,4).Value = "completed"

# Row 49
$ws.Cells.Item(49,1).Value = 45706
$ws.Cells.Item(49,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(49,2).Value = "Day by day class"

# Row 50
$ws.Cells.Item(50,1).Value = 45707
$ws.Cells.Item(50,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(50,2).Value = "Assignment12(merges)"
$ws.Cells.Item(50,4).Value = "completed"

# Row 51
$ws.Cells.Item(51,1).Value = 45708
$ws.Cells.Item(51,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(51,2).Value = "no class"

# Row 52
$ws.Cells.Item(52,1).Value = 45709
$ws.Cells.Item(52,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(52,2).Value = "leave"

# Row 53
$ws.Cells.Item(53,1).Value = 45710
$ws.Cells.Item(53,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(53,2).Value = "Saturday"

# Row 54
$ws.Cells.Item(54,1).Value = 45711
$ws.Cells.Item(54,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(54,2).Value = "Sunday"

# Row 55
$ws.Cells.Item(55,1).Value = 45712
$ws.Cells.Item(55,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(55,2).Value = "leave"

# Row 56
$ws.Cells.Item(56,1).Value = 45713
$ws.Cells.Item(56,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(56,2).Value = "leave"

# Row 57
$ws.Cells.Item(57,1).Value = 45714
$ws.Cells.Item(57,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(57,2).Value = "leave"

# Row 58
$ws.Cells.Item(58,1).Value = 45715
$ws.Cells.Item(58,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(58,2).Value = "leave"

# Row 59
$ws.Cells.Item(59,1).Value = 45716
$ws.Cells.Item(59,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(59,2).Value = "Assignment12(merges)"
$ws.Cells.Item(59,4).Value = "completed"

# Row 60
$ws.Cells.Item(60,1).Value = 45717
$ws.Cells.Item(60,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(60,2).Value = "Saturday"

# Row 61
$ws.Cells.Item(61,1).Value = 45718
$ws.Cells.Item(61,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(61,2).Value = "Sunday"

# Row 62
$ws.Cells.Item(62,1).Value = 45719
$ws.Cells.Item(62,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(62,2).Value = "leave"

# Row 63
$ws.Cells.Item(63,1).Value = 45720
$ws.Cells.Item(63,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(63,2).Value = "leave"

# Row 64
$ws.Cells.Item(64,1).Value = 45721
$ws.Cells.Item(64,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(64,2).Value = "leave"

# Row 65
$ws.Cells.Item(65,1).Value = 45722
$ws.Cells.Item(65,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(65,2).Value = "leave"

# Row 66
$ws.Cells.Item(66,1).Value = 45723
$ws.Cells.Item(66,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(66,2).Value = "leave"

# Row 67
$ws.Cells.Item(67,1).Value = 45724
$ws.Cells.Item(67,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(67,2).Value = "Saturday"

# Row 68
$ws.Cells.Item(68,1).Value = 45725
$ws.Cells.Item(68,1).NumberFormat = "m/d/yy"
$ws.Cells.Item(68,2).Value = "Sunday"

# Move the selection to match where the author last edited (D68)
$ws.Range("D68").Select() | Out-Null